# Remove the word "maximal" from the k-core / k-dense definitions on
# slide 3 and slide 4 ("rimossa parole massimale da slide 3 e 4").
#
# Slide 3 ("K-core"): "... is a set of maximal connected sub-graphs ..."
#                  -> "... is a set of connected sub-graphs ..."
# Slide 4 ("K-dense"): "... is a set of maximal connected sub-graphs ..."
#                  -> "... is a set of connected sub-graphs ..."

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 3 - "K-core"
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange

$full3 = $tr3.Text
$needle3 = "maximal connected sub-"
$idx3 = $full3.IndexOf($needle3)

# First collapse "maximal connected sub-" down to "connected sub-"
# (keeps it as its own run, separate from the preceding "a set of " run).
$run3b = $tr3.Characters($idx3 + 1, $needle3.Length)
$run3b.Text = "connected sub-"

# Now shrink the preceding "a set of " run down to "a set " and move the
# "of " back onto the following run, yielding the three runs:
#   "a set "  /  "of connected "  /  "sub-"
$full3b = $tr3.Text
$needleA3 = "a set of "
$idxA3 = $full3b.IndexOf($needleA3)
$runA3 = $tr3.Characters($idxA3 + 1, $needleA3.Length)
$runA3.Text = "a set "

$startB3 = $idxA3 + 1 + "a set ".Length
$runB3 = $tr3.Characters($startB3, "connected sub-".Length)
[void]$runB3.InsertBefore("of ")

# Finally split "of connected sub-" into "of connected " and "sub-" so the
# run boundaries match the target (sub- keeps its own run).
$full3c = $tr3.Text
$startSub3 = $idxA3 + 1 + "a set ".Length + "of connected ".Length
$runSub3 = $tr3.Characters($startSub3, "sub-".Length)
$runSub3.Text = "sub-"

# ---------------------------------------------------------------
# Slide 4 - "K-dense"
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange

$full4 = $tr4.Text
$needleB4 = "connected sub-"
$idxB4 = $full4.IndexOf($needleB4)
$runB4 = $tr4.Characters($idxB4 + 1, $needleB4.Length)
$runB4.Text = "sub-"

$full4b = $tr4.Text
$needleA4 = "set of maximal "
$idxA4 = $full4b.IndexOf($needleA4)
$runA4 = $tr4.Characters($idxA4 + 1, $needleA4.Length)
$runA4.Text = "set of connected "
